# Auto-generated Excel COM-interop script reproducing the workbook diff.
# For each sheet, update the changed numeric cells (H:N columns) to match
# the target state described by the commit diff. A few cells are newly
# added (set where previously blank) and one cell (ALC!N101) is cleared
# entirely since the target state has no value there.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 67 cell updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 442
$ws.Range("I6").Value = 442
$ws.Range("K6").Value = 1326
$ws.Range("M6").Value = -1214
$ws.Range("H8").Value = 53.42857
$ws.Range("I8").Value = 53.42857
$ws.Range("K8").Value = 160.28571
$ws.Range("M8").Value = -21.28570999999999
$ws.Range("H58").Value = 2271
$ws.Range("J58").Value = 4250
$ws.Range("L58").Value = 12750
$ws.Range("N58").Value = -13050
$ws.Range("H62").Value = 2118.8
$ws.Range("I62").Value = 2174.5
$ws.Range("K62").Value = 2174.5
$ws.Range("M62").Value = -1550.5
$ws.Range("H65").Value = 2118.8
$ws.Range("I65").Value = 2174.5
$ws.Range("K65").Value = 10872.5
$ws.Range("M65").Value = -7752.5
$ws.Range("H96").Value = 381.6
$ws.Range("I96").Value = 377
$ws.Range("K96").Value = 1131
$ws.Range("M96").Value = 242
$ws.Range("H99").Value = 700.8182
$ws.Range("I99").Value = 611.2857
$ws.Range("J99").Value = 857.5
$ws.Range("K99").Value = 1833.8571
$ws.Range("L99").Value = 2572.5
$ws.Range("M99").Value = -335.8571000000002
$ws.Range("N99").Value = -5568.5
$ws.Range("H100").Value = 4001.4666
$ws.Range("I100").Value = 4361.091
$ws.Range("J100").Value = 3012.5
$ws.Range("K100").Value = 4361.091
$ws.Range("L100").Value = 3012.5
$ws.Range("M100").Value = -3820.091
$ws.Range("N100").Value = -4094.5
$ws.Range("H101").Value = 393.6
$ws.Range("I101").Value = 393.6
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1180.8
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 441.1999999999998
$ws.Range("N101").ClearContents()
$ws.Range("H113").Value = 5394.1055
$ws.Range("I113").Value = 4166.5
$ws.Range("J113").Value = 5960.6924
$ws.Range("K113").Value = 4166.5
$ws.Range("L113").Value = 5960.6924
$ws.Range("M113").Value = -912.5
$ws.Range("N113").Value = -12468.6924
$ws.Range("H116").Value = 4500
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -11884
$ws.Range("H132").Value = 983.625
$ws.Range("I132").Value = 981.2857
$ws.Range("K132").Value = 2943.8571
$ws.Range("M132").Value = -413.8571000000002
$ws.Range("H137").Value = 2893.4119
$ws.Range("I137").Value = 2015.6666
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 6046.9998
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -3496.9998
$ws.Range("N137").Value = -20100

# --- Sheet ARM: 20 cell updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1749.75
$ws.Range("I45").Value = 1666.3334
$ws.Range("K45").Value = 1666.3334
$ws.Range("M45").Value = -1289.3334
$ws.Range("H63").Value = 4003
$ws.Range("I63").Value = 3000
$ws.Range("K63").Value = 3000
$ws.Range("M63").Value = -2314
$ws.Range("H66").Value = 4003
$ws.Range("I66").Value = 3000
$ws.Range("K66").Value = 15000
$ws.Range("M66").Value = -11568
$ws.Range("H97").Value = 593
$ws.Range("J97").Value = 595
$ws.Range("L97").Value = 595
$ws.Range("N97").Value = -1587
$ws.Range("H112").Value = 39499.5
$ws.Range("J112").Value = 39499.5
$ws.Range("L112").Value = 39499.5
$ws.Range("N112").Value = -42453.5

# --- Sheet BSM: 27 cell updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4541
$ws.Range("I86").Value = 3644.2856
$ws.Range("K86").Value = 3644.2856
$ws.Range("M86").Value = -2521.2856
$ws.Range("H89").Value = 4541
$ws.Range("I89").Value = 3644.2856
$ws.Range("K89").Value = 18221.428
$ws.Range("M89").Value = -12605.428
$ws.Range("H94").Value = 1677.8235
$ws.Range("I94").Value = 1460.1538
$ws.Range("K94").Value = 1460.1538
$ws.Range("M94").Value = -1009.1538
$ws.Range("H105").Value = 2920.375
$ws.Range("I105").Value = 2844.8333
$ws.Range("K105").Value = 2844.8333
$ws.Range("M105").Value = -1097.8333
$ws.Range("H110").Value = 148752
$ws.Range("J110").Value = 148752
$ws.Range("L110").Value = 148752
$ws.Range("N110").Value = -156932
$ws.Range("H134").Value = 6097.9
$ws.Range("I134").Value = 6612.5
$ws.Range("J134").Value = 4039.5
$ws.Range("K134").Value = 19837.5
$ws.Range("L134").Value = 12118.5
$ws.Range("M134").Value = -17302.5
$ws.Range("N134").Value = -17188.5

# --- Sheet CRP: 29 cell updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2398.1667
$ws.Range("I99").Value = 2217.8
$ws.Range("J99").Value = 3300
$ws.Range("K99").Value = 2217.8
$ws.Range("L99").Value = 3300
$ws.Range("M99").Value = -719.8000000000002
$ws.Range("N99").Value = -6296
$ws.Range("H105").Value = 889.8333
$ws.Range("I105").Value = 835
$ws.Range("K105").Value = 835
$ws.Range("M105").Value = 912
$ws.Range("H126").Value = 2398.1667
$ws.Range("I126").Value = 2217.8
$ws.Range("J126").Value = 3300
$ws.Range("K126").Value = 6653.400000000001
$ws.Range("L126").Value = 9900
$ws.Range("M126").Value = -4183.400000000001
$ws.Range("N126").Value = -14840
$ws.Range("H132").Value = 3568.3
$ws.Range("J132").Value = 3465.25
$ws.Range("L132").Value = 10395.75
$ws.Range("N132").Value = -15455.75
$ws.Range("H134").Value = 3879.875
$ws.Range("I134").Value = 3590.6667
$ws.Range("J134").Value = 4747.5
$ws.Range("K134").Value = 10772.0001
$ws.Range("L134").Value = 14242.5
$ws.Range("M134").Value = -8237.000100000001
$ws.Range("N134").Value = -19312.5

# --- Sheet CUL: 34 cell updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1602.7916
$ws.Range("I4").Value = 1459.2778
$ws.Range("K4").Value = 4377.8334
$ws.Range("M4").Value = -4265.8334
$ws.Range("H52").Value = 1487
$ws.Range("J52").Value = 1487
$ws.Range("L52").Value = 4461
$ws.Range("N52").Value = -4993
$ws.Range("H68").Value = 1176
$ws.Range("I68").Value = 1074.5
$ws.Range("J68").Value = 1205
$ws.Range("K68").Value = 3223.5
$ws.Range("L68").Value = 3615
$ws.Range("M68").Value = -2412.5
$ws.Range("N68").Value = -5237
$ws.Range("H71").Value = 1176
$ws.Range("I71").Value = 1074.5
$ws.Range("J71").Value = 1205
$ws.Range("K71").Value = 9670.5
$ws.Range("L71").Value = 10845
$ws.Range("M71").Value = -5614.5
$ws.Range("N71").Value = -18957
$ws.Range("H108").Value = 584.4
$ws.Range("I108").Value = 584.4
$ws.Range("K108").Value = 1753.2
$ws.Range("M108").Value = 1126.8
$ws.Range("H139").Value = 1640.4546
$ws.Range("I139").Value = 672.44446
$ws.Range("K139").Value = 2017.33338
$ws.Range("M139").Value = 3122.66662
$ws.Range("H140").Value = 3890.5
$ws.Range("I140").Value = 3509.6667
$ws.Range("K140").Value = 10529.0001
$ws.Range("M140").Value = -5349.000100000001

# --- Sheet GSM: 8 cell updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 12131.167
$ws.Range("I102").Value = 3223.5454
$ws.Range("K102").Value = 3223.5454
$ws.Range("M102").Value = -1601.5454
$ws.Range("H122").Value = 3289.3
$ws.Range("J122").Value = 3398
$ws.Range("L122").Value = 10194
$ws.Range("N122").Value = -15094

# --- Sheet LTW: 2 cell updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1500
$ws.Range("H27").Value = 1500

# --- Sheet WVR: 22 cell updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 999
$ws.Range("I81").Value = 999
$ws.Range("K81").Value = 1998
$ws.Range("M81").Value = -937
$ws.Range("H84").Value = 999
$ws.Range("I84").Value = 999
$ws.Range("K84").Value = 9990
$ws.Range("M84").Value = -4686
$ws.Range("H113").Value = 334.83334
$ws.Range("I113").Value = 336
$ws.Range("J113").Value = 332.5
$ws.Range("K113").Value = 1008
$ws.Range("L113").Value = 997.5
$ws.Range("M113").Value = 1162
$ws.Range("N113").Value = -5337.5
$ws.Range("H136").Value = 7802.6665
$ws.Range("I136").Value = 3406.8333
$ws.Range("J136").Value = 10733.223
$ws.Range("K136").Value = 10220.4999
$ws.Range("L136").Value = 32199.669
$ws.Range("M136").Value = -7670.499899999999
$ws.Range("N136").Value = -37299.669
